# Voto servicio comprobado y funcion modificarDatos
# Adds row 22 to Hoja1: "11 - Puntuar" entry documenting the new voting route.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row content (mirrors the style of rows 16-21: D/F/G/H/J populated).
# Values are written in the order the new route/endpoint text first, so new
# shared-string entries land in the same order as the source workbook.
$ws.Range("G22").Value = 'app.post("/servicios/votar/:id_servicio/:id_solucionador",isUser, voteServicio);'
$ws.Range("F22").Value = "Método (POST)"
$ws.Range("D22").Value = "11 - Puntuar"
$ws.Range("H22").Value = "isUser"
$ws.Range("J22").Value = "ok"

# G22 gets a dedicated new style: same font/fill as the rest of the table
# rows, but with wrap text enabled (the endpoint text is long).
$ws.Range("D16").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").WrapText = $true

# Column G widened to fit the new, longer route string; wrap text means it's
# no longer a "best fit" width.
$ws.Columns.Item(7).ColumnWidth = 74.33

# Update selection to mirror where the author left off after adding the row
$ws.Range("H23").Select()
